$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "charrange" column (F) entries for rows 11-18,
# mirroring the existing byterange (D) / charindex (E) columns.
$ws.Range("F11").Value = "get_bytevector_from_charrange"
$ws.Range("F12").Value = "get_charvector_from_charrange"
$ws.Range("F13").Value = "get_glyphvector_from_charrange"
$ws.Range("F14").Value = "get_byteiterator_from_charrange"
$ws.Range("F15").Value = "get_chariterator_from_charrange"
$ws.Range("F16").Value = "get_glyphiterator_from_charrange"
$ws.Range("F17").Value = "get_strref_from_charrange"
$ws.Range("F18").Value = "get_string_from_charrange"

# F11:F14 pick up the wrap-text / top-vertical alignment style already
# used by the analogous cells in column D (byterange).
$ws.Range("F11:F14").WrapText = $true
$ws.Range("F11:F18").VerticalAlignment = -4160

# Move the active selection from F10 to F9, as in the final sheet state.
$ws.Range("F9").Select()
